$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark D-column cells we are about to rewrite as Text so Excel
# does not reinterpret numeric-looking strings (e.g. "0.999") as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Step 2: write the updated cell values
$ws.Range("D2").Value = "69.685.65"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "3.550.37"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "198.52"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "587.26"
$ws.Range("E6").Value = "  -3.07%  "
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.208"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "0.632"
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").Value = "52.36"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").Value = "0.0000291"
$ws.Range("E12").Value = "  -4.92%  "
$ws.Range("D13").Value = "689.78"
$ws.Range("E13").Value = "  +16.13%  "
$ws.Range("D14").Value = "9.36"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "4.107.22"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").Value = "69.707.63"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "3.552.31"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "12.52"
$ws.Range("E18").Value = "  -5.77%  "
$ws.Range("D19").Value = "18.67"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "0.975"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").Value = "18.40"
$ws.Range("E22").Value = "  +3.80%  "
$ws.Range("D23").Value = "108.78"
$ws.Range("E23").Value = "  +5.98%  "
$ws.Range("D24").Value = "5.24"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("E25").Value = "  -4.41%  "
$ws.Range("D26").Value = "2.98"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").Value = "10.31"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("D28").Value = "9.80"
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("D29").Value = "33.74"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "4.41"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").Value = "6.94"
$ws.Range("D32").Value = "11.97"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("E33").Value = "  -3.92%  "
$ws.Range("D34").Value = "62.16"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").Value = "3.808.31"
$ws.Range("E35").Value = "  -3.67%  "
$ws.Range("D36").Value = "0.0₃0824"
$ws.Range("E36").Value = "  -8.16%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "3.70"
$ws.Range("E38").Value = "  +4.15%  "
$ws.Range("E39").Value = "  -6.96%  "
$ws.Range("D40").Value = "499.21"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("E41").Value = "  -4.53%  "
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").Value = "34.95"
$ws.Range("E43").Value = "  -6.35%  "
$ws.Range("D44").Value = "0.0463"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("D45").Value = "2.95"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "8.45"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "1.83"
$ws.Range("E50").Value = "  +21.85%  "
$ws.Range("D51").Value = "2.80"
$ws.Range("E51").Value = "  +69.11%  "

# Step 3: restore default cell style on the D-column cells touched above
# (keeps them text-typed without leaving a lingering explicit format).
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
